$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @("Paris Saint-Germain", 4.790490341753343, 10.05987055016181, 0.7806379971734303, 12, 0.04882364506193066, 42)
    3  = @("RC Lens", 5.05668449197861, 7.456171735241503, 0.5144346431435445, 5, 0.07494466111521028, 38)
    5  = @("Olympique de Marseille", 6.049226441631505, 11.71895424836601, 0.6540903540903541, 14, 0.05218113834648941, 5)
    8  = @("RC Strasbourg", 6.67248322147651, 10.43594306049822, 0.4993581514762516, 21, 0.04277673545966229, 19)
    9  = @("Toulouse", 6.726723095525998, 6.019891500904159, 0.4153869545385434, 17, 0.08081077476996933, -10)
    10 = @("Angers", 6.544242424242424, 8.064631956912029, 0.3580377159435579, 13, 0.06854049072687478, -14)
    11 = @("AS Monaco", 4.396590909090909, 7.390109890109891, 0.5751889848812095, 17, 0.0705264221481257, 22)
    12 = @("Lorient", 6.607194244604316, 8.461215932914046, 0.3828875767048744, 11, 0.07785804380231624, -9)
    14 = @("Le Havre", 4.909199522102748, 8.03671706263499, 0.4069129209316125, 13, 0.09330323551542513, -4)
    15 = @("Nice", 5.396662387676509, 7.94589552238806, 0.4510786360473208, 22, 0.07235745248416138, -7)
    16 = @("Paris FC", 5.967651195499297, 11.94143167028199, 0.5210970464135021, 21, 0.07520420617782368, -10)
    17 = @("Auxerre", 5.901960784313726, 6.142329020332717, 0.4652948847993045, 14, 0.0901985111662531, -23)
    18 = @("Nantes", 5.768308921438082, 7.452290076335878, 0.363308812490617, 17, 0.09586776859504133, -25)
}

foreach ($r in $data.Keys) {
    $row = $data[$r]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
}
